# Apply "repull data" update to the dSF column (F) for the winckowski_josh
# workbook: a handful of rows were recalculated/repulled and their dSF
# value changed. Update just those cells to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = -2
    5  = -1
    8  = -1
    12 = 5
    13 = -4
    15 = 2
    18 = -5
    22 = 0
    25 = 0
    31 = 0
    32 = 5
    33 = 1
    34 = -1
    35 = -1
    38 = 3
    40 = -2
    41 = 2
    43 = -2
    46 = 4
    47 = 1
    51 = 0
    53 = 7
    54 = 2
    55 = -1
    56 = 2
    59 = -3
    63 = 0
    64 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
